# MirrorMe - modifications made to accommodate for the refactoring changes
#
# For both worksheets ("MirroMe Voorbeeldredenering" and "MirrorMe Example Argument"):
#  - G2 changes from "[Template,]" / "[Template,]" to "Template" (matches A2/G1 wording)
#  - A new column H is introduced, mirroring column G in rows 1 and 2 (the "requires" /
#    "Template" header block), and splitting the combined "X, Y" values that used to live
#    in G5 and G7 into separate G/H cells.

$wb = $excel.ActiveWorkbook

$sheetNames = @("MirroMe Voorbeeldredenering", "MirrorMe Example Argument")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: the "[Template,]" label becomes plain "Template"
    $ws.Range("G2").Value2 = "Template"

    # New column H mirrors column G for the header rows 1 and 2
    $ws.Range("H1").Value2 = $ws.Range("G1").Value2
    $ws.Range("H2").Value2 = $ws.Range("G2").Value2

    # Row 5: split the combined "A, B" value across G5/H5
    $parts5 = $ws.Range("G5").Value2 -split ", "
    $ws.Range("G5").Value2 = $parts5[0]
    $ws.Range("H5").Value2 = $parts5[1]

    # Row 7: split the combined "A, B" value across G7/H7
    $parts7 = $ws.Range("G7").Value2 -split ", "
    $ws.Range("G7").Value2 = $parts7[0]
    $ws.Range("H7").Value2 = $parts7[1]
}

# Row 16 on the "MirrorMe Example Argument" sheet had its explicit row height reset
# (it becomes an auto-height row) as part of this refactor
$ws2Row16 = $wb.Worksheets.Item("MirrorMe Example Argument")
$ws2Row16.Rows("16:16").AutoFit()

# Restore/normalize selection on both sheets to H7, keep sheet1 as the active tab
$ws2 = $wb.Worksheets.Item("MirrorMe Example Argument")
$ws2.Activate()
$ws2.Range("H7").Select()

$ws1 = $wb.Worksheets.Item("MirroMe Voorbeeldredenering")
$ws1.Activate()
$ws1.Range("H7").Select()
